# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.006876353814593728, 0.3375848360084654, 0.7127328510149897, 246.9852506941017, 0, 248.0424447349397),
    @(0.3464964993005633, 3486408257.158922, 3.082599426703578, 71517.89157740913, 1, 3486479778.479595),
    @(0.7287194209349384, 9.226618575922256, 3.082599426703578, 246.9852506941017, 0, 260.0231881176624),
    @(1.505614041169197, 1766.335244827366, 3.082599426703578, 71517.89157740913, 1, 73288.81503570436),
    @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 0, 10.35301142835362),
    @(0.006876353814593728, 0.004309184025731883, 3.082599426703578, 6.48142807727062, 0, 9.575213041814525),
    @(0.000000000546653389221774, 1.65323645889881, 2938.103010863317, 5548678842208.939, 0, 5548678845148.695)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

$wb.Save()
